# Sync automatico del tracker (cada 3h)
# Appends the latest scraped rows (event_id, fecha, jugador_A, jugador_B,
# pronostico, cuota, resultado, profit) to the bottom of the results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to sync. event_id / fecha are kept as text (matching how the
# scraper emits them), resultado/profit are still pending so they land
# as empty text cells, and cuota is numeric.
$newRows = @(
    @{ Row = 205; EventId = "14316422"; Fecha = "2025-08-15"; JugadorA = "Varvara Gracheva";     JugadorB = "Veronika Kudermetova";       Pronostico = "Gana Varvara Gracheva";          Cuota = 2.63 },
    @{ Row = 206; EventId = "14368955"; Fecha = "2025-08-15"; JugadorA = "Timofey Skatov";        JugadorB = "Valentin Vacherot";           Pronostico = "Gana Timofey Skatov";            Cuota = 2.5  },
    @{ Row = 207; EventId = "14368951"; Fecha = "2025-08-15"; JugadorA = "Stefano Travaglia";     JugadorB = "Juan Carlos Prado Angelo";    Pronostico = "Gana Juan Carlos Prado Angelo";  Cuota = 2.25 },
    @{ Row = 208; EventId = "14437028"; Fecha = "2025-08-15"; JugadorA = "Maximilian Neuchrist";  JugadorB = "Lilian Marmousez";            Pronostico = "Gana Maximilian Neuchrist";      Cuota = 2.62 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Columns A (event_id) and B (fecha) look numeric/date-like to Excel's
    # smart parsing, so force them to Text first, write, then reset the
    # cell style back to Normal so no stray style index sticks around.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $r.EventId
    $ws.Cells.Item($row, 1).Style = "Normal"

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $r.Fecha
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = $r.JugadorA
    $ws.Cells.Item($row, 4).Value = $r.JugadorB
    $ws.Cells.Item($row, 5).Value = $r.Pronostico

    $ws.Cells.Item($row, 6).Value = $r.Cuota

    # resultado / profit are still pending (no result yet) -> empty text cells.
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 7).Style = "Normal"

    $ws.Cells.Item($row, 8).Value = "'"
    $ws.Cells.Item($row, 8).Style = "Normal"
}
